# Autogenerated on Mon Feb 09 2015 03:30:35 GMT+0000 (Coordinated Universal Time)
#
# Inserts the MSME size-classification table (Number of employees / Assets /
# Turnover, broken out by Micro / Small / Medium / Large) into rows 18-22 of
# the Cambodia Summary sheet, and relocates the existing "NIS" source footer
# (previously rows 23-24) down to rows 29-30 to make room for it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Pull the existing footer text (rows 23-24) off the sheet before we
#        move anything, then delete those two rows outright so nothing is
#        left behind at the old location. ---
$sourceName = $ws.Range("A23").Text
$sourceDesc = $ws.Range("A24").Text
$ws.Rows("23:24").Delete()

# --- 2. New table header (row 18), bold like the other section titles. ---
$ws.Range("B18").Value = "Number of employees"
$ws.Range("C18").Value = "Assets (local currency, unless noted otherwise)"
$ws.Range("D18").Value = "Turnover (local currency, unless noted otherwise)"
$ws.Range("B18:D18").Font.Bold = $true

# --- 3. Table body (rows 19-22), plain/default formatting. ---
$ws.Range("A19").Value = "Micro"
$ws.Range("B19").Value = "<=10"
$ws.Range("C19").Value = "< USD50,000"
$ws.Range("D19").Value = ""

$ws.Range("A20").Value = "Small"
$ws.Range("B20").Value = "11-50"
$ws.Range("C20").Value = "USD50,000 - USD 250,000"
$ws.Range("D20").Value = ""

$ws.Range("A21").Value = "Medium"
$ws.Range("B21").Value = "51-100"
$ws.Range("C21").Value = "USD 250,000 - USD 500,000"
$ws.Range("D21").Value = ""

$ws.Range("A22").Value = "Large"
$ws.Range("B22").Value = ">100"
$ws.Range("C22").Value = ">USD 500,000"
$ws.Range("D22").Value = ""

# --- 4. Re-create the footer two rows further down (rows 29-30), preserving
#        the original "title" (bold) / "source" (italic) look. ---
$ws.Range("A29").Value = $sourceName
$ws.Range("A29").Font.Bold = $true

$ws.Range("A30").Value = $sourceDesc
$ws.Range("A30").Font.Italic = $true
